$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 1028.8334
$ws.Range("I31").Value = 1028.8334
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 3086.5002
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -2856.5002

$ws.Range("H33").Value = 302.81818
$ws.Range("I33").Value = 227.05882
$ws.Range("J33").Value = 560.4
$ws.Range("K33").Value = 227.05882
$ws.Range("L33").Value = 560.4
$ws.Range("M33").Value = 1.941180000000003
$ws.Range("N33").Value = -1018.4

$ws.Range("H64").Value = 3227.389
$ws.Range("I64").Value = 3011.625
$ws.Range("J64").Value = 3400
$ws.Range("K64").Value = 3011.625
$ws.Range("L64").Value = 3400
$ws.Range("M64").Value = -2763.625
$ws.Range("N64").Value = -3896

$ws.Range("H67").Value = 3227.389
$ws.Range("I67").Value = 3011.625
$ws.Range("J67").Value = 3400
$ws.Range("K67").Value = 3011.625
$ws.Range("L67").Value = 3400
$ws.Range("M67").Value = -2153.625
$ws.Range("N67").Value = -5116

$ws.Range("H86").Value = 21698.8
$ws.Range("I86").Value = 100000
$ws.Range("J86").Value = 2123.5
$ws.Range("K86").Value = 100000
$ws.Range("L86").Value = 2123.5
$ws.Range("M86").Value = -98877
$ws.Range("N86").Value = -4369.5

$ws.Range("H89").Value = 21698.8
$ws.Range("I89").Value = 100000
$ws.Range("J89").Value = 2123.5
$ws.Range("K89").Value = 500000
$ws.Range("L89").Value = 10617.5
$ws.Range("M89").Value = -494384
$ws.Range("N89").Value = -21849.5

$ws.Range("H134").Value = 49343
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 49343
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 49343
$ws.Range("N134").Value = -59483

$ws.Range("H137").Value = 1193.6666
$ws.Range("I137").Value = 789.73914
$ws.Range("J137").Value = 1682.6316
$ws.Range("K137").Value = 2369.21742
$ws.Range("L137").Value = 5047.8948
$ws.Range("M137").Value = 180.7825800000001
$ws.Range("N137").Value = -10147.8948

$ws.Range("H138").Value = 1647.37
$ws.Range("I138").Value = 825.0417
$ws.Range("J138").Value = 2406.4424
$ws.Range("K138").Value = 2475.1251
$ws.Range("L138").Value = 7219.3272
$ws.Range("M138").Value = 2664.8749
$ws.Range("N138").Value = -17499.3272

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 553.6667
$ws.Range("I14").Value = 553.6667
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 553.6667
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -378.6667

$ws.Range("H22").Value = 600
$ws.Range("I22").Value = 600
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 600
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -301

$ws.Range("H61").Value = 1623.5
$ws.Range("I61").Value = 1525.3334
$ws.Range("J61").Value = 2507
$ws.Range("K61").Value = 1525.3334
$ws.Range("L61").Value = 2507
$ws.Range("M61").Value = -1313.3334
$ws.Range("N61").Value = -2931

$ws.Range("H132").Value = 21742458
$ws.Range("I132").Value = 27028024
$ws.Range("J132").Value = 12914
$ws.Range("K132").Value = 81084072
$ws.Range("L132").Value = 38742
$ws.Range("M132").Value = -81081542
$ws.Range("N132").Value = -43802

$ws.Range("H136").Value = 1623.5
$ws.Range("I136").Value = 1525.3334
$ws.Range("J136").Value = 2507
$ws.Range("K136").Value = 4576.0002
$ws.Range("L136").Value = 7521
$ws.Range("M136").Value = -2026.0002
$ws.Range("N136").Value = -12621

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1712308.5
$ws.Range("I134").Value = 1025.6136
$ws.Range("J134").Value = 5297853.5
$ws.Range("K134").Value = 3076.8408
$ws.Range("L134").Value = 15893560.5
$ws.Range("M134").Value = -541.8407999999999
$ws.Range("N134").Value = -15898630.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4004.6667
$ws.Range("I62").Value = 2850.6667
$ws.Range("J62").Value = 4581.6665
$ws.Range("K62").Value = 2850.6667
$ws.Range("L62").Value = 4581.6665
$ws.Range("M62").Value = -2226.6667
$ws.Range("N62").Value = -5829.6665

$ws.Range("H65").Value = 4004.6667
$ws.Range("I65").Value = 2850.6667
$ws.Range("J65").Value = 4581.6665
$ws.Range("K65").Value = 14253.3335
$ws.Range("L65").Value = 22908.3325
$ws.Range("M65").Value = -11133.3335
$ws.Range("N65").Value = -29148.3325

$ws.Range("H99").Value = 58832296
$ws.Range("I99").Value = 100012390
$ws.Range("J99").Value = 3584.2856
$ws.Range("K99").Value = 100012390
$ws.Range("L99").Value = 3584.2856
$ws.Range("M99").Value = -100010892
$ws.Range("N99").Value = -6580.2856

$ws.Range("H126").Value = 58832296
$ws.Range("I126").Value = 100012390
$ws.Range("J126").Value = 3584.2856
$ws.Range("K126").Value = 300037170
$ws.Range("L126").Value = 10752.8568
$ws.Range("M126").Value = -300034700
$ws.Range("N126").Value = -15692.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 45477796
$ws.Range("I44").Value = 83335250
$ws.Range("J44").Value = 48854.8
$ws.Range("K44").Value = 250005750
$ws.Range("L44").Value = 146564.4
$ws.Range("M44").Value = -250005352
$ws.Range("N44").Value = -147360.4

$ws.Range("H46").Value = 560.8
$ws.Range("I46").Value = 200
$ws.Range("J46").Value = 2004
$ws.Range("K46").Value = 600
$ws.Range("L46").Value = 6012
$ws.Range("M46").Value = -509
$ws.Range("N46").Value = -6194

$ws.Range("H137").Value = 35716360
$ws.Range("I137").Value = 38463616
$ws.Range("J137").Value = 2000
$ws.Range("K137").Value = 115390848
$ws.Range("L137").Value = 6000
$ws.Range("M137").Value = -115385748
$ws.Range("N137").Value = -16200

$ws.Range("H139").Value = 296537.88
$ws.Range("I139").Value = 1778.6666
$ws.Range("J139").Value = 1003960
$ws.Range("K139").Value = 5335.9998
$ws.Range("L139").Value = 3011880
$ws.Range("M139").Value = -195.9997999999996
$ws.Range("N139").Value = -3022160

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").ClearContents()

$ws.Range("H132").Value = 5319.1787
$ws.Range("I132").Value = 881.7
$ws.Range("J132").Value = 16412.875
$ws.Range("K132").Value = 2645.1
$ws.Range("L132").Value = 49238.625
$ws.Range("M132").Value = -115.1000000000004
$ws.Range("N132").Value = -54298.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 921.1316
$ws.Range("I93").Value = 900.1111
$ws.Range("J93").Value = 972.7273
$ws.Range("K93").Value = 900.1111
$ws.Range("L93").Value = 972.7273
$ws.Range("M93").Value = 347.8889
$ws.Range("N93").Value = -3468.7273

$ws.Range("H132").Value = 19237476
$ws.Range("I132").Value = 28572966
$ws.Range("J132").Value = 17347.53
$ws.Range("K132").Value = 85718898
$ws.Range("L132").Value = 52042.59
$ws.Range("M132").Value = -85716368
$ws.Range("N132").Value = -57102.59

$ws.Range("H136").Value = 4125.75
$ws.Range("I136").Value = 3800.8823
$ws.Range("J136").Value = 5966.6665
$ws.Range("K136").Value = 11402.6469
$ws.Range("L136").Value = 17899.9995
$ws.Range("M136").Value = -8852.6469
$ws.Range("N136").Value = -22999.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 22519.717
$ws.Range("I132").Value = 22011.164
$ws.Range("J132").Value = 28749.5
$ws.Range("K132").Value = 66033.492
$ws.Range("L132").Value = 86248.5
$ws.Range("M132").Value = -63503.492
$ws.Range("N132").Value = -91308.5

$ws.Range("H136").Value = 6564.478
$ws.Range("I136").Value = 9790.666999999999
$ws.Range("J136").Value = 3045
$ws.Range("K136").Value = 29372.001
$ws.Range("L136").Value = 9135
$ws.Range("M136").Value = -26822.001
$ws.Range("N136").Value = -14235

$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws.Range("H139").Value = 48771.285
$ws.Range("I139").Value = 45004.5
$ws.Range("J139").Value = 50278
$ws.Range("K139").Value = 45004.5
$ws.Range("L139").Value = 50278
$ws.Range("M139").Value = -39864.5
$ws.Range("N139").Value = -60558
